$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1, Q1 following the pattern of row 1 (0-based sequence), with same style as O1
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I column -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K column -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M column -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O column -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P column = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q column = 2
}
